$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. "总计" sheet: insert a new row for 2022-Q3 at the top of the
#    data (row 2), push existing quarters down, renumber index col
# ---------------------------------------------------------------
$ws1.Rows.Item(2).Insert()

# Re-apply the index-column style (bold/border/centered) that got
# shifted down to A3 onto the freshly inserted A2 cell
$ws1.Range("A3").Copy($ws1.Range("A2"))

# The row-insert carries the headers bold/border style onto the
# new rows other cells -- reset B2:D2 back to the plain/default style
$ws1.Range("B2:D2").Style = "Normal"

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 13
$ws1.Range("D2").Value = 0.33

# Renumber the 0-based index column for the rows pushed down
$ws1.Range("A3").Value = 1
$ws1.Range("A4").Value = 2
$ws1.Range("A5").Value = 3
$ws1.Range("A6").Value = 4
$ws1.Range("A7").Value = 5

# ---------------------------------------------------------------
# 2. Add a new "2022-Q3" worksheet, positioned right after "总计"
#    (i.e. before the existing "2022-Q2" sheet)
# ---------------------------------------------------------------
$sheetBefore = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($sheetBefore)
$newSheet.Name = "2022-Q3"

# Match the outline / page-margin conventions used by the other
# quarterly sheets in this workbook
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$ps = $newSheet.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# Header row (B1:H1): copy the bold/border/centered header style
# used throughout the workbook, then set the header captions
foreach ($col in @("B","C","D","E","F","G","H")) {
    $ws1.Range("B1").Copy($newSheet.Range($col + "1"))
}
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows 2-14: column A (0-based index) and column H (rank)
# are numeric; columns B-G are stored as text (fund codes must
# keep leading zeros, and the decimal figures must keep their
# fixed-precision formatting, e.g. "0.30" rather than 0.3)
$data = @(
    @(0, '320022', '诺安研究精选股票', '6.17', '92.67', '1.97', '0.1215', 10),
    @(1, '001097', '华泰柏瑞积极优选股票A', '1.20', '83.01', '3.15', '0.0378', 4),
    @(2, '014839', '兴银碳中和主题混合C', '0.67', '92.10', '4.96', '0.0332', 5),
    @(3, '010797', '长城优选回报六个月持有期混合A', '2.83', '28.94', '1.05', '0.0297', 7),
    @(4, '000066', '诺安鸿鑫混合A', '0.69', '79.09', '4.14', '0.0286', 4),
    @(5, '014838', '兴银碳中和主题混合A', '0.55', '92.10', '4.96', '0.0273', 5),
    @(6, '014831', '兴银中证1000指数增强A', '1.37', '83.33', '1.23', '0.0169', 2),
    @(7, '009937', '东方欣益一年持有期偏债混合A', '2.22', '20.97', '0.72', '0.0160', 8),
    @(8, '014832', '兴银中证1000指数增强C', '0.90', '83.33', '1.23', '0.0111', 2),
    @(9, '016283', '华泰柏瑞积极优选股票C', '0.16', '83.01', '3.15', '0.0050', 4),
    @(10, '010798', '长城优选回报六个月持有期混合C', '0.30', '28.94', '1.05', '0.0032', 7),
    @(11, '009938', '东方欣益一年持有期偏债混合C', '0.37', '20.97', '0.72', '0.0027', 8),
    @(12, '014498', '诺安鸿鑫混合C', '0.01', '79.09', '4.14', '0.0004', 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $idxCell = $newSheet.Range("A" + $r)
    $ws1.Range("A2").Copy($idxCell)
    $idxCell.Value = $row[0]

    foreach ($col in @("B","C","D","E","F","G")) {
        $newSheet.Range($col + $r).NumberFormat = "@"
    }
    $newSheet.Range("B" + $r).Value = $row[1]
    $newSheet.Range("C" + $r).Value = $row[2]
    $newSheet.Range("D" + $r).Value = $row[3]
    $newSheet.Range("E" + $r).Value = $row[4]
    $newSheet.Range("F" + $r).Value = $row[5]
    $newSheet.Range("G" + $r).Value = $row[6]
    $newSheet.Range("B" + $r + ":G" + $r).Style = "Normal"

    $newSheet.Range("H" + $r).Value = $row[7]
}

# Keep "总计" as the active sheet/tab, same as before the edit
$ws1.Activate()
$ws1.Range("A1").Select() | Out-Null

Write-Host "Edit complete"
